# Update the "取得日時" (acquired datetime) timestamps in rows 2-11 of the
# "ランサーズ" sheet from 2025-09-27 18:23:03 to 2025-09-27 18:28:10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-27 18:28:10"

for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
